# Update TPM-derived NATMI ligand-receptor metrics (Sema3f -> Nrp2) after rerunning scripts with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value2 = 40.83537666666667
$ws.Range("H2").Value2 = 122.50613
$ws.Range("I2").Value2 = 0.9274830900091532
$ws.Range("J2").Value2 = 0.9274830900091531
$ws.Range("M2").Value2 = 57.65261933333333
$ws.Range("N2").Value2 = 172.957858
$ws.Range("O2").Value2 = 0.6817060950001529
$ws.Range("P2").Value2 = 0.6817060950001529
$ws.Range("Q2").Value2 = 2354.266426296615
$ws.Range("R2").Value2 = 21188.39783666954
$ws.Range("S2").Value2 = 0.6322708754688151
$ws.Range("T2").Value2 = 0.6322708754688151

# Row 3
$ws.Range("G3").Value2 = 40.83537666666667
$ws.Range("H3").Value2 = 122.50613
$ws.Range("I3").Value2 = 0.9274830900091532
$ws.Range("J3").Value2 = 0.9274830900091531
$ws.Range("O3").Value2 = 0.1019529789289588
$ws.Range("P3").Value2 = 0.1019529789289588
$ws.Range("Q3").Value2 = 352.0937792896222
$ws.Range("R3").Value2 = 3168.8440136066
$ws.Range("S3").Value2 = 0.09455966393266882
$ws.Range("T3").Value2 = 0.0945596639326688

# Row 4
$ws.Range("G4").Value2 = 40.83537666666667
$ws.Range("H4").Value2 = 122.50613
$ws.Range("I4").Value2 = 0.9274830900091532
$ws.Range("J4").Value2 = 0.9274830900091531
$ws.Range("M4").Value2 = 2.790736
$ws.Range("N4").Value2 = 8.372208
$ws.Range("O4").Value2 = 0.0329987043561157
$ws.Range("P4").Value2 = 0.0329987043561157
$ws.Range("Q4").Value2 = 113.9607557372267
$ws.Range("R4").Value2 = 1025.64680163504
$ws.Range("S4").Value2 = 0.0306057402825087
$ws.Range("T4").Value2 = 0.03060574028250869

# Row 5
$ws.Range("G5").Value2 = 40.83537666666667
$ws.Range("H5").Value2 = 122.50613
$ws.Range("I5").Value2 = 0.9274830900091532
$ws.Range("J5").Value2 = 0.9274830900091531
$ws.Range("M5").Value2 = 15.50544933333333
$ws.Range("N5").Value2 = 46.516348
$ws.Range("O5").Value2 = 0.1833422217147727
$ws.Range("P5").Value2 = 0.1833422217147727
$ws.Range("Q5").Value2 = 633.1708639125823
$ws.Range("R5").Value2 = 5698.537775213241
$ws.Range("S5").Value2 = 0.1700468103251606
$ws.Range("T5").Value2 = 0.1700468103251606

# Row 6
$ws.Range("I6").Value2 = 0.03813623414934058
$ws.Range("J6").Value2 = 0.03813623414934057
$ws.Range("M6").Value2 = 57.65261933333333
$ws.Range("N6").Value2 = 172.957858
$ws.Range("O6").Value2 = 0.6817060950001529
$ws.Range("P6").Value2 = 0.6817060950001529
$ws.Range("Q6").Value2 = 96.80268745632111
$ws.Range("R6").Value2 = 871.22418710689
$ws.Range("S6").Value2 = 0.02599770325995844
$ws.Range("T6").Value2 = 0.02599770325995844

# Row 7
$ws.Range("I7").Value2 = 0.03813623414934058
$ws.Range("J7").Value2 = 0.03813623414934057
$ws.Range("O7").Value2 = 0.1019529789289588
$ws.Range("P7").Value2 = 0.1019529789289588
$ws.Range("S7").Value2 = 0.00388810267665756
$ws.Range("T7").Value2 = 0.003888102676657559

# Row 8
$ws.Range("I8").Value2 = 0.03813623414934058
$ws.Range("J8").Value2 = 0.03813623414934057
$ws.Range("M8").Value2 = 2.790736
$ws.Range("N8").Value2 = 8.372208
$ws.Range("O8").Value2 = 0.0329987043561157
$ws.Range("P8").Value2 = 0.0329987043561157
$ws.Range("Q8").Value2 = 4.685836444293334
$ws.Range("R8").Value2 = 42.17252799864001
$ws.Range("S8").Value2 = 0.001258446315949693
$ws.Range("T8").Value2 = 0.001258446315949693

# Row 9
$ws.Range("I9").Value2 = 0.03813623414934058
$ws.Range("J9").Value2 = 0.03813623414934057
$ws.Range("M9").Value2 = 15.50544933333333
$ws.Range("N9").Value2 = 46.516348
$ws.Range("O9").Value2 = 0.1833422217147727
$ws.Range("P9").Value2 = 0.1833422217147727
$ws.Range("Q9").Value2 = 26.03470896970444
$ws.Range("R9").Value2 = 234.31238072734
$ws.Range("S9").Value2 = 0.006991981896774886
$ws.Range("T9").Value2 = 0.006991981896774885

# Row 10
$ws.Range("G10").Value2 = 1.503819
$ws.Range("H10").Value2 = 4.511457
$ws.Range("I10").Value2 = 0.03415584247746153
$ws.Range("J10").Value2 = 0.03415584247746152
$ws.Range("M10").Value2 = 57.65261933333333
$ws.Range("N10").Value2 = 172.957858
$ws.Range("O10").Value2 = 0.6817060950001529
$ws.Range("P10").Value2 = 0.6817060950001529
$ws.Range("Q10").Value2 = 86.69910435323399
$ws.Range("R10").Value2 = 780.291939179106
$ws.Range("S10").Value2 = 0.02328424599675065
$ws.Range("T10").Value2 = 0.02328424599675064

# Row 11
$ws.Range("G11").Value2 = 1.503819
$ws.Range("H11").Value2 = 4.511457
$ws.Range("I11").Value2 = 0.03415584247746153
$ws.Range("J11").Value2 = 0.03415584247746152
$ws.Range("O11").Value2 = 0.1019529789289588
$ws.Range("P11").Value2 = 0.1019529789289588
$ws.Range("Q11").Value2 = 12.96633846186
$ws.Range("R11").Value2 = 116.69704615674
$ws.Range("S11").Value2 = 0.003482289888405472
$ws.Range("T11").Value2 = 0.003482289888405471

# Row 12
$ws.Range("G12").Value2 = 1.503819
$ws.Range("H12").Value2 = 4.511457
$ws.Range("I12").Value2 = 0.03415584247746153
$ws.Range("J12").Value2 = 0.03415584247746152
$ws.Range("M12").Value2 = 2.790736
$ws.Range("N12").Value2 = 8.372208
$ws.Range("O12").Value2 = 0.0329987043561157
$ws.Range("P12").Value2 = 0.0329987043561157
$ws.Range("Q12").Value2 = 4.196761820784
$ws.Range("R12").Value2 = 37.770856387056
$ws.Range("S12").Value2 = 0.001127098547947811
$ws.Range("T12").Value2 = 0.001127098547947811

# Row 13
$ws.Range("G13").Value2 = 1.503819
$ws.Range("H13").Value2 = 4.511457
$ws.Range("I13").Value2 = 0.03415584247746153
$ws.Range("J13").Value2 = 0.03415584247746152
$ws.Range("M13").Value2 = 15.50544933333333
$ws.Range("N13").Value2 = 46.516348
$ws.Range("O13").Value2 = 0.1833422217147727
$ws.Range("P13").Value2 = 0.1833422217147727
$ws.Range("Q13").Value2 = 23.317389311004
$ws.Range("R13").Value2 = 209.856503799036
$ws.Range("S13").Value2 = 0.006262208044357602
$ws.Range("T13").Value2 = 0.006262208044357601

# Row 14
$ws.Range("E14").Value2 = 1
$ws.Range("F14").Value2 = 0.3333333333333333
$ws.Range("G14").Value2 = 0.009899
$ws.Range("H14").Value2 = 0.029697
$ws.Range("I14").Value2 = 0.0002248333640447365
$ws.Range("J14").Value2 = 0.0002248333640447365
$ws.Range("M14").Value2 = 57.65261933333333
$ws.Range("N14").Value2 = 172.957858
$ws.Range("O14").Value2 = 0.6817060950001529
$ws.Range("P14").Value2 = 0.6817060950001529
$ws.Range("Q14").Value2 = 0.5707032787806666
$ws.Range("R14").Value2 = 5.136329509026
$ws.Range("S14").Value2 = 0.0001532702746286851
$ws.Range("T14").Value2 = 0.0001532702746286851

# Row 15
$ws.Range("E15").Value2 = 1
$ws.Range("F15").Value2 = 0.3333333333333333
$ws.Range("G15").Value2 = 0.009899
$ws.Range("H15").Value2 = 0.029697
$ws.Range("I15").Value2 = 0.0002248333640447365
$ws.Range("J15").Value2 = 0.0002248333640447365
$ws.Range("O15").Value2 = 0.1019529789289588
$ws.Range("P15").Value2 = 0.1019529789289588
$ws.Range("Q15").Value2 = 0.08535188372666666
$ws.Range("R15").Value2 = 0.76816695354
$ws.Range("S15").Value2 = 0.00002292243122697995
$ws.Range("T15").Value2 = 0.00002292243122697995

# Row 16
$ws.Range("E16").Value2 = 1
$ws.Range("F16").Value2 = 0.3333333333333333
$ws.Range("G16").Value2 = 0.009899
$ws.Range("H16").Value2 = 0.029697
$ws.Range("I16").Value2 = 0.0002248333640447365
$ws.Range("J16").Value2 = 0.0002248333640447365
$ws.Range("M16").Value2 = 2.790736
$ws.Range("N16").Value2 = 8.372208
$ws.Range("O16").Value2 = 0.0329987043561157
$ws.Range("P16").Value2 = 0.0329987043561157
$ws.Range("Q16").Value2 = 0.027625495664
$ws.Range("R16").Value2 = 0.248629460976
$ws.Range("S16").Value2 = 0.000007419209709503196
$ws.Range("T16").Value2 = 0.000007419209709503193

# Row 17
$ws.Range("E17").Value2 = 1
$ws.Range("F17").Value2 = 0.3333333333333333
$ws.Range("G17").Value2 = 0.009899
$ws.Range("H17").Value2 = 0.029697
$ws.Range("I17").Value2 = 0.0002248333640447365
$ws.Range("J17").Value2 = 0.0002248333640447365
$ws.Range("M17").Value2 = 15.50544933333333
$ws.Range("N17").Value2 = 46.516348
$ws.Range("O17").Value2 = 0.1833422217147727
$ws.Range("P17").Value2 = 0.1833422217147727
$ws.Range("Q17").Value2 = 0.1534884429506667
$ws.Range("R17").Value2 = 1.381395986556
$ws.Range("S17").Value2 = 0.00004122144847956829
$ws.Range("T17").Value2 = 0.00004122144847956829
